# Apply updated crypto prices/volume percentages (row 2-51: columns D "Price", E "Volume(1h)")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (D) text. Rows not listed here keep their original Price value.
$priceUpdates = @{
    2 = "42.313.88"
    3 = "2.274.33"
    4 = "1.00"
    5 = "306.04"
    6 = "97.36"
    9 = "0.493"
    10 = "35.88"
    11 = "0.0798"
    14 = "2.601.34"
    15 = "14.47"
    16 = "2.276.37"
    18 = "42.215.24"
    19 = "12.52"
    21 = "6.02"
    22 = "67.91"
    23 = "239.02"
    28 = "37.43"
    29 = "9.57"
    30 = "2.11"
    31 = "159.60"
    33 = "0.999"
    34 = "3.20"
    35 = "0.0742"
    36 = "17.36"
    43 = "1.989.70"
    44 = "0.0286"
    45 = "19.00"
    46 = "2.94"
    47 = "9.95"
    48 = "53.35"
    49 = "1.52"
    50 = "72.23"
    51 = "91.83"
}

# Map of row -> new Volume(1h) text (E column).
$volumeUpdates = @{
    2 = "  +1.20%  "
    3 = "  +0.45%  "
    4 = "  +0.02%  "
    5 = "  +0.95%  "
    6 = "  +5.27%  "
    7 = "  +0.18%  "
    8 = "  +0.02%  "
    9 = "  +1.96%  "
    10 = "  +10.81%  "
    11 = "  +0.13%  "
    12 = "  -0.76%  "
    13 = "  +0.09%  "
    14 = "  -0.85%  "
    15 = "  +1.40%  "
    16 = "  +0.24%  "
    17 = "  +2.19%  "
    18 = "  +1.09%  "
    19 = "  -1.73%  "
    20 = "  +0.40%  "
    21 = "  +1.43%  "
    22 = "  +0.88%  "
    23 = "  -1.79%  "
    24 = "  +0.50%  "
    25 = "  +1.29%  "
    26 = "  -0.06%  "
    27 = "  -0.71%  "
    28 = "  +7.11%  "
    29 = "  +0.24%  "
    30 = "  +1.97%  "
    31 = "  -0.22%  "
    32 = "  +0.68%  "
    33 = "  +0.03%  "
    34 = "  +6.14%  "
    35 = "  +0.16%  "
    36 = "  +3.08%  "
    37 = "  +0.31%  "
    38 = "  -0.60%  "
    39 = "  +2.80%  "
    40 = "  -0.98%  "
    41 = "  +4.55%  "
    42 = "  +14.40%  "
    43 = "  -1.37%  "
    44 = "  +1.23%  "
    45 = "  -3.23%  "
    46 = "  +1.85%  "
    47 = "  -3.99%  "
    48 = "  +0.69%  "
    49 = "  +1.18%  "
    50 = "  -1.21%  "
    51 = "  +0.20%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage so numeric-looking strings (e.g. "1.00") keep their exact formatting
    # instead of being auto-converted to numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

# Restore default (General) cell style now that the text values are safely stored,
# so the cell formatting matches the original (unstyled) cells.
foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Style = "Normal"
}
